# ---------------------------------------------------------------------------
# Applies the RFI-template edit:
#   1. Row 1: drop the empty B1/C1 placeholder cells (keep A1:C1 merged).
#   2. Backup & Disaster Recovery section: add a new question row
#      "Is there any backup equipment on site?" right after the existing
#      backup questions, pushing the "Documentation & Handoff" section
#      (and its two questions) down by one row.
# ---------------------------------------------------------------------------

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A scratch cell used to stash formatting via copy/paste-special; it lives far
# outside the used range so it never disturbs real data.
$scratch = "Z500"

# ---------------------------------------------------------------------------
# 1) Row 1 - remove the blank, styled B1/C1 cells while preserving the
#    merged A1:C1 title cell (style, value and row height).
# ---------------------------------------------------------------------------

$titleValue  = $ws.Range("A1").Value()
$titleHeight = $ws.Rows.Item(1).RowHeight()

# Stash A1's current format (style 1) so we can re-apply it after re-merging.
$ws.Range("A1").Copy()
$ws.Range($scratch).PasteSpecial(-4122)   # xlPasteFormats

# Unmerge, then reset the whole row to the default "Normal" style/contents.
# Merging while every cell is "Normal" means Excel does not stamp a
# non-default style onto B1/C1, so they are omitted from the saved XML.
$ws.Range("A1:C1").UnMerge()
$ws.Range("A1:C1").ClearContents()
$ws.Range("A1:C1").Style = "Normal"
$ws.Range("A1:C1").Merge()

# Re-apply the title's format/value/row height now that B1/C1 stay untouched.
$ws.Range($scratch).Copy()
$ws.Range("A1").PasteSpecial(-4122)       # xlPasteFormats
$ws.Range("A1").Value = $titleValue
$ws.Rows.Item(1).RowHeight = $titleHeight

$ws.Range($scratch).Delete()

# ---------------------------------------------------------------------------
# 2) Insert the new "Is there any backup equipment on site?" question.
#
#    Before:                              After:
#    103 (ht=24) Documentation & Handoff  103 (ht=24) Is there any backup...
#    104         Is there existing IT...  104         Documentation & Handoff
#    105         What admin credentials.  105         Is there existing IT...
#                                          106         What admin credentials.
# ---------------------------------------------------------------------------

# Insert a fresh blank row at 104. Row 103 (with its ht=24 custom height and
# the "Documentation & Handoff" header) is left exactly where it is; the old
# rows 104/105 shift down to 105/106.
$ws.Rows.Item(104).Insert()

# Give the new row104 a default style before merging (same trick as above)
# so that merging A104:C104 does not create phantom B104/C104 cells.
$ws.Range("A104").Style = "Normal"
$ws.Range("A104:C104").Merge()

# Move the section header down into row 104: copy A103's current format
# (style 3, the section-header look) onto A104, then set its text.
$ws.Range("A103").Copy()
$ws.Range("A104").PasteSpecial(-4122)     # xlPasteFormats
$ws.Range("A104").Value = "Documentation & Handoff"

# Free up row 103: remove its old header merge, then turn it into a normal
# question row matching the style of the other Backup & DR questions.
$ws.Range("A103:C103").UnMerge()
$ws.Range("A102:C102").Copy()
$ws.Range("A103:C103").PasteSpecial(-4122) # xlPasteFormats (s=4 / s=5)
$ws.Range("A103").Value = "Is there any backup equipment on site?"

# Rows 105/106 already contain the right questions (shifted automatically by
# the Insert above), so nothing else needs to change.
